$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 and C2 are stored as text in the workbook (e.g. "42", "3"); use a
# leading apostrophe so Excel keeps the new values ("13", "5") as text
# instead of auto-converting them to numbers.
$ws.Range("B2").Value = "'13"
$ws.Range("C2").Value = "'5"

# D2 and E2 are plain numeric cells.
$ws.Range("D2").Value = 0.1
$ws.Range("E2").Value = 0.04
